$wb = $excel.ActiveWorkbook

# --- Sheet "Personnes" (sheet1): remove the last row (idPers=3, NOM, Prenom, 0) ---
$wsPersonnes = $wb.Worksheets.Item("Personnes")
$wsPersonnes.Rows.Item(4).Delete()

# --- Sheet "Adresses" (sheet2): update row 4 and add a new row 5 ---
$wsAdresses = $wb.Worksheets.Item("Adresses")

$wsAdresses.Range("A4").Value = 3
$wsAdresses.Range("B4").Value = 10
$wsAdresses.Range("C4").Value = "rue"
$wsAdresses.Range("D4").Value = 86000
$wsAdresses.Range("E4").Value = "ville"

$wsAdresses.Range("A5").Value = 4
$wsAdresses.Range("B5").Value = 10
$wsAdresses.Range("C5").Value = "newrue"
$wsAdresses.Range("D5").Value = 87000
$wsAdresses.Range("E5").Value = "newville"
